# Updates market/profit figures on several leve sheets, as pulled by the
# scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# ALC sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 70
$ws.Cells.Item(70, 8).Value  = 1350       # H70
$ws.Cells.Item(70, 9).Value  = 970        # I70
$ws.Cells.Item(70, 11).Value = 2910       # K70
$ws.Cells.Item(70, 13).Value = -2640      # M70

# Row 73
$ws.Cells.Item(73, 8).Value  = 1350       # H73
$ws.Cells.Item(73, 9).Value  = 970        # I73
$ws.Cells.Item(73, 11).Value = 2910       # K73
$ws.Cells.Item(73, 13).Value = -1974      # M73

# Rows 125-141: the pricing columns (H:N) are no longer populated for
# these entries, so clear them out entirely.
$ws.Range("H125:N141").ClearContents()

# ----------------------------------------------------------------------
# ARM sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 35
$ws.Cells.Item(35, 8).Value  = 0          # H35
$ws.Cells.Item(35, 9).Value  = 0          # I35
$ws.Cells.Item(35, 11).Value = 0          # K35
$ws.Range("M35").ClearContents()          # M35 removed

# ----------------------------------------------------------------------
# BSM sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 29
$ws.Cells.Item(29, 8).Value  = 1458       # H29
$ws.Cells.Item(29, 9).Value  = 516        # I29
$ws.Cells.Item(29, 10).Value = 2400       # J29
$ws.Cells.Item(29, 11).Value = 516        # K29
$ws.Cells.Item(29, 12).Value = 2400       # L29
$ws.Cells.Item(29, 13).Value = -227       # M29
$ws.Cells.Item(29, 14).Value = -2978      # N29 (new)

# Row 30
$ws.Cells.Item(30, 8).Value  = 6011       # H30
$ws.Cells.Item(30, 9).Value  = 0          # I30
$ws.Cells.Item(30, 10).Value = 6011       # J30
$ws.Cells.Item(30, 11).Value = 0          # K30
$ws.Cells.Item(30, 12).Value = 6011       # L30
$ws.Range("M30").ClearContents()          # M30 removed
$ws.Cells.Item(30, 14).Value = -6261      # N30 (new)

# Row 36
$ws.Cells.Item(36, 8).Value  = 2000       # H36
$ws.Cells.Item(36, 9).Value  = 2000       # I36
$ws.Cells.Item(36, 11).Value = 2000       # K36
$ws.Cells.Item(36, 13).Value = -1466      # M36

# Row 37
$ws.Cells.Item(37, 8).Value  = 1000       # H37
$ws.Cells.Item(37, 9).Value  = 1000       # I37
$ws.Cells.Item(37, 11).Value = 1000       # K37
$ws.Cells.Item(37, 13).Value = -863       # M37

# ----------------------------------------------------------------------
# CRP sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 81
$ws.Cells.Item(81, 8).Value  = 48576.8    # H81
$ws.Cells.Item(81, 10).Value = 48576.8    # J81
$ws.Cells.Item(81, 12).Value = 48576.8    # L81
$ws.Cells.Item(81, 14).Value = -50572.8   # N81

# Row 84
$ws.Cells.Item(84, 8).Value  = 48576.8    # H84
$ws.Cells.Item(84, 10).Value = 48576.8    # J84
$ws.Cells.Item(84, 12).Value = 145730.4   # L84
$ws.Cells.Item(84, 14).Value = -155714.4  # N84

# ----------------------------------------------------------------------
# CUL sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 102
$ws.Cells.Item(102, 8).Value  = 14000      # H102
$ws.Cells.Item(102, 10).Value = 14000      # J102
$ws.Cells.Item(102, 12).Value = 42000      # L102
$ws.Cells.Item(102, 14).Value = -46868     # N102

# Row 123
$ws.Cells.Item(123, 8).Value  = 6416.6665  # H123
$ws.Cells.Item(123, 9).Value  = 0          # I123
$ws.Cells.Item(123, 10).Value = 6416.6665  # J123
$ws.Cells.Item(123, 11).Value = 0          # K123
$ws.Cells.Item(123, 12).Value = 19249.9995 # L123
$ws.Range("M123").ClearContents()          # M123 removed
$ws.Cells.Item(123, 14).Value = -24149.9995 # N123

# ----------------------------------------------------------------------
# LTW sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 139
$ws.Cells.Item(139, 8).Value  = 57950     # H139
$ws.Cells.Item(139, 10).Value = 57950     # J139
$ws.Cells.Item(139, 12).Value = 57950     # L139
$ws.Cells.Item(139, 14).Value = -68230    # N139

# ----------------------------------------------------------------------
# WVR sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 107
$ws.Cells.Item(107, 8).Value  = 996.26666  # H107
$ws.Cells.Item(107, 9).Value  = 1239.909   # I107
$ws.Cells.Item(107, 10).Value = 326.25     # J107
$ws.Cells.Item(107, 11).Value = 3719.727   # K107
$ws.Cells.Item(107, 12).Value = 978.75     # L107
$ws.Cells.Item(107, 13).Value = -1799.727  # M107
$ws.Cells.Item(107, 14).Value = -4818.75   # N107
